$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overhead")

# ---------------------------------------------------------------------------
# 1. Fill in the new TRIM rows (4-11) with PN / description data.
#    The order in which new string values are first written controls the
#    order new entries are appended to the shared-strings table, so the
#    writes below are intentionally sequenced to match the source order.
# ---------------------------------------------------------------------------

$ws.Range("B4").Value = "TRIM"

$ws.Range("C4").Value  = "832Z6501-1"
$ws.Range("C8").Value  = "832Z6501-111"
$ws.Range("C5").Value  = "832Z6501-2"
$ws.Range("C9").Value  = "832Z6501-112"
$ws.Range("C6").Value  = "832Z6501-3"
$ws.Range("C10").Value = "832Z6501-113"
$ws.Range("C7").Value  = "832Z6501-4"
$ws.Range("C11").Value = "832Z6501-114"

$ws.Range("E4").Value  = "CÓ LỖ - GẮN VỚI SCREW PN BACS12ER3K7 (24INCH - 35LB VỚI CTR, 50LB VỚI OUTB)"
$ws.Range("E8").Value  = "KHÔNG LỖ (24INCH - 35LB VỚI CTR, 50LB VỚI OUTB)"
$ws.Range("E5").Value  = "CÓ LỖ - GẮN VỚI SCREW PN BACS12ER3K7 (36INCH - 56LB VỚI CTR, 80LB VỚI OUTB)"
$ws.Range("E9").Value  = "KHÔNG LỖ (36INCH - 56LB VỚI CTR, 80LB VỚI OUTB)"
$ws.Range("E6").Value  = "CÓ LỖ - GẮN VỚI SCREW PN BACS12ER3K7 (42INCH - 66LB VỚI CTR, 95LB VỚI OUTB)"
$ws.Range("E10").Value = "KHÔNG LỖ (42INCH - 66LB VỚI CTR, 95LB VỚI OUTB)"
$ws.Range("E7").Value  = "CÓ LỖ - GẮN VỚI SCREW PN BACS12ER3K7 (48INCH - 77LB VỚI CTR, 111LB VỚI OUTB)"
$ws.Range("E11").Value = "KHÔNG LỖ (48INCH - 77LB VỚI CTR, 111LB VỚI OUTB)"

# Remaining TRIM cells in column B just repeat the same label.
$ws.Range("B5").Value = "TRIM"
$ws.Range("B6").Value = "TRIM"
$ws.Range("B7").Value = "TRIM"
$ws.Range("B8").Value = "TRIM"
$ws.Range("B9").Value = "TRIM"
$ws.Range("B10").Value = "TRIM"
$ws.Range("B11").Value = "TRIM"

# Rows 4-11 now wrap onto two lines, matching row height bump.
$ws.Range("A4:E11").RowHeight = 30

# ---------------------------------------------------------------------------
# 2. Re-bucket the A/C column by aircraft type now that four new B787 TRIM
#    rows exist: rows 4-13 => B787, 14-17 => ATR72, 18-23 => A321 (new).
# ---------------------------------------------------------------------------

$ws.Range("A8:A13").Value = "B787"
$ws.Range("A14:A17").Value = "ATR72"

# Extend the table down to row 23, copying the formatting already used for
# the existing blank A321 rows, then label the new rows A321.
$ws.Range("A17:E17").Copy() | Out-Null
$ws.Range("A18:E23").PasteSpecial(-4122) | Out-Null
$ws.Range("A18:A23").Value = "A321"

# Column E got a bit wider to accommodate the new long descriptions.
$ws.Columns.Item(5).ColumnWidth = 51.7

# ---------------------------------------------------------------------------
# 3. Restore focus back onto the Overhead sheet (it had moved to Pax seat).
# ---------------------------------------------------------------------------

$ws.Activate()
$ws.Range("E16").Select() | Out-Null
